# Update Milestone 3 - Directional, Point and Spot light working -
# Instantiating properly done - Input layout changed

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mark milestone completion (Column F = "X") for several rows
$ws.Range("F4").Value = "X"
$ws.Range("F5").Value = "X"
$ws.Range("F18").Value = "X"
$ws.Range("F23").Value = "X"
$ws.Range("F30").Value = "X"
$ws.Range("F31").Value = "X"
$ws.Range("F32").Value = "X"
$ws.Range("F65").Value = "X"
$ws.Range("F66").Value = "X"

# Row 6 also gets its Student milestone (Column E) set to "II"
$ws.Range("E6").Value = "II"
$ws.Range("F6").Value = "X"

# Move the active selection on Sheet1 to F18 (reflecting the user's last edit)
[void]$ws.Range("F18").Select()
